$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.003.76"
$ws.Range("D2").ClearFormats()

$ws.Range("D3").Value = "'2.287.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +4.97%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'252.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.46%  "

$ws.Range("D6").Value = "'0.631"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.00%  "

$ws.Range("D7").Value = "'72.85"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +9.51%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.664"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +12.70%  "

$ws.Range("D10").Value = "'39.45"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.55%  "

$ws.Range("D11").Value = "'0.0981"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.98%  "

$ws.Range("D12").Value = "'59.84"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.55%  "

$ws.Range("E13").Value = "  +8.97%  "

$ws.Range("D14").Value = "'0.106"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").Value = "'2.628.66"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.98%  "

$ws.Range("D16").Value = "'15.13"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.91%  "

$ws.Range("D17").Value = "'0.894"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.77%  "

$ws.Range("D18").Value = "'2.283.32"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.55%  "

$ws.Range("D19").Value = "'42.940.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.47%  "

$ws.Range("E20").Value = "  +6.64%  "

$ws.Range("D21").Value = "'6.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.94%  "

$ws.Range("D22").Value = "'73.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.96%  "

$ws.Range("D23").Value = "'238.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.55%  "

$ws.Range("D24").Value = "'2.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.84%  "

$ws.Range("D25").Value = "'3.90"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.77%  "

$ws.Range("D26").Value = "'11.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.38%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").Value = "'2.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.58%  "

$ws.Range("D29").Value = "'3.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").Value = "'2.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +8.26%  "

$ws.Range("D31").Value = "'168.19"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").Value = "'21.21"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.10%  "

$ws.Range("D33").Value = "'6.34"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +11.71%  "

$ws.Range("D34").Value = "'0.128"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.07%  "

$ws.Range("D35").Value = "'0.0811"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.67%  "

$ws.Range("D36").Value = "'31.36"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +26.95%  "

$ws.Range("D37").Value = "'4.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +21.59%  "

$ws.Range("D38").Value = "'0.127"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.81%  "

$ws.Range("D39").Value = "'4.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.26%  "

$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("E41").Value = "  +6.20%  "

$ws.Range("D42").Value = "'13.41"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +18.81%  "

$ws.Range("D43").Value = "'6.11"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +10.49%  "

$ws.Range("E44").Value = "  +13.13%  "

$ws.Range("D45").Value = "'9.22"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.26%  "

$ws.Range("E46").Value = "  -9.75%  "

$ws.Range("D47").Value = "'61.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("E48").Value = "  +4.95%  "

$ws.Range("E49").Value = "  +5.09%  "

$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.20"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.01%  "
